# Update "南宁-漫展信息.xlsx" to add the "南宁·2024三月三国潮动漫节（良牙春典）" event
# to the "展览" (sheet 1) and "全部类型" (sheet 4) sheets, just above the last
# ("南宁·布谷鸟动漫展4th") row, and bump a couple of "想去人数" counters.
#
# NOTE: this interop layer's function calls only reliably marshal COM object
# arguments when passed *positionally* - named arguments (e.g. "-ws $ws1")
# silently drop the binding, so every call below uses positional args.

$wb = $excel.ActiveWorkbook

# xlPasteFormats constant used to copy the bordered/bold "index column" style.
$xlPasteFormats = -4122

function Add-SanyuesanEvent {
    param(
        $ws,
        [int]$strawberryRow,      # row with "南宁·草莓动漫节" (想去人数 1736 -> 1740)
        [int]$daCgRow,            # row with "南宁·第一届ANE·DACG动漫嘉年华" (769 -> 773)
        [int]$insertRow,          # row number the new event row should occupy
        [int]$oldIndexValue,      # value that belongs in column A of the new row
        [int]$newLastIndexValue,  # value that belongs in column A of the shifted last row
        [int]$newLastCount        # updated "想去人数" value for the shifted last row
    )

    # Bump the two "想去人数" counters that changed independently of the new row.
    $ws.Cells.Item($strawberryRow, 6).Value = 1740
    $ws.Cells.Item($daCgRow, 6).Value = 773

    # Insert a brand-new row directly above the current last row, pushing the
    # last row (南宁·布谷鸟动漫展4th) one row further down.
    $ws.Rows.Item($insertRow).Insert()

    # Match the formatting (border/bold/centering) used by the other index cells.
    $ws.Cells.Item($insertRow - 1, 1).Copy() | Out-Null
    $ws.Cells.Item($insertRow, 1).PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false

    $ws.Cells.Item($insertRow, 1).Value = $oldIndexValue
    $ws.Cells.Item($insertRow, 2).Value = "'2024-05-01"
    $ws.Cells.Item($insertRow, 3).Value = "南宁·2024三月三国潮动漫节（良牙春典）"
    $ws.Cells.Item($insertRow, 4).Value = "民族大道106号 南宁国际会展中心"
    $ws.Cells.Item($insertRow, 5).Value = "2024.05.01 09:30-05.02 17:30"
    $ws.Cells.Item($insertRow, 6).Value = 34
    $ws.Cells.Item($insertRow, 7).Value = 55
    $ws.Cells.Item($insertRow, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82416"
    $ws.Cells.Item($insertRow, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/Df29DEWI1709708373277.jpeg"

    # Update the index + count of the row that just got shifted down.
    $shiftedRow = $insertRow + 1
    $ws.Cells.Item($shiftedRow, 1).Value = $newLastIndexValue
    $ws.Cells.Item($shiftedRow, 6).Value = $newLastCount
}

# Sheet 1 = "展览": rows 4/5 are 草莓动漫节/ANE·DACG; last row was row 6 -> new
# row 6 inserted (old row 6, 布谷鸟动漫展4th, becomes row 7).
$ws1 = $wb.Worksheets.Item(1)
Add-SanyuesanEvent $ws1 4 5 6 5 6 196

# Sheet 4 = "全部类型": row 4 is 草莓动漫节, row 6 is ANE·DACG (row 5 is the
# concert); last row was row 7 -> new row 7 inserted (old row 7, 布谷鸟动漫展4th,
# becomes row 8).
$ws4 = $wb.Worksheets.Item(4)
Add-SanyuesanEvent $ws4 4 6 7 6 7 196

Write-Host "edit complete"
